$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 220, shifting existing rows 220:285 down to 221:286
$ws.Rows.Item(220).Insert()

# Populate the newly inserted row 220 with the new record
$ws.Range("A220").Value = 6
$ws.Range("B220").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C220").Value = "Metropolitana"
$ws.Range("D220").Value = 44785
$ws.Range("E220").Value = 13
$ws.Range("F220").Value = 100112026
$ws.Range("G220").Value = "Haba"
$ws.Range("H220").Value = "Sin especificar"
$ws.Range("I220").Value = "Primera"
$ws.Range("J220").Value = 400
$ws.Range("K220").Value = 17000
$ws.Range("L220").Value = 18000
$ws.Range("M220").Value = 17425
$ws.Range("N220").Value = '$/saco 25 kilos'
$ws.Range("O220").Value = "Región de Coquimbo"
$ws.Range("P220").Value = 697
$ws.Range("Q220").Value = 25
$ws.Range("R220").Value = "Hortaliza"
